# 3.4.2 Suicide mortality rate - add a new "2023" data column (T) to the table
# and refresh the layout (column widths / row heights) to match the
# published sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -----------------------------------------------------
# Columns A:C become one uniform width, column D gets narrower.
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 30.666666666666668
$ws.Range("D1").EntireColumn.ColumnWidth = 8.833333333333334

# --- Row heights ---------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 16.5

# --- New column T: copy formatting from column S (same row), then set
#     the 2023 values -----------------------------------------------------
$ws.Range("S4").Copy($ws.Range("T4"))
$ws.Range("T4").Value = 2023

$ws.Range("S5").Copy($ws.Range("T5"))
$ws.Range("T5").Value = 4.8

$ws.Range("S6").Copy($ws.Range("T6"))
$ws.Range("T6").Value = 5.7

$ws.Range("S7").Copy($ws.Range("T7"))
$ws.Range("T7").Value = 1.9

$ws.Range("S8").Copy($ws.Range("T8"))
$ws.Range("T8").Value = 8.9

$ws.Range("S9").Copy($ws.Range("T9"))
$ws.Range("T9").Value = 11.9

$ws.Range("S10").Copy($ws.Range("T10"))
$ws.Range("T10").Value = 2.5

$ws.Range("S11").Copy($ws.Range("T11"))
$ws.Range("T11").Value = 0.7

$ws.Range("S12").Copy($ws.Range("T12"))
$ws.Range("T12").Value = 12.7

$ws.Range("S13").Copy($ws.Range("T13"))
$ws.Range("T13").Value = 1.1

$ws.Range("S14").Copy($ws.Range("T14"))
$ws.Range("T14").Value = 2.2

# --- Reset the saved selection back to the top-left cell ----------------
$ws.Range("A1").Select()
